$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.178.72'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '3.818.32'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '702.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '3.817.30'
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.52'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").Value = '4.463.51'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").Value = '3.820.24'
$ws.Range("E16").Value = '  -4.30%  '
$ws.Range("D17").Value = '71.221.14'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '510.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").Value = '3.970.97'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  -4.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.173'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("D37").Value = '3.781.44'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '171.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '427.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.31%  '
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.33%  '
